$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the existing FY-2013..FY-2022 index labels from rows 2-11 first, so the shared-string
# table forgets their old slot order; we re-add all labels below in chronological order so
# FY-2003..FY-2022 end up contiguous, matching how the source data export was regenerated.
$ws.Range("A2:D11").ClearContents()

# Row labels (column A) for all data rows, now extended back to FY-2003, written oldest-to-newest
$ws.Range("A2").Value = "NVDA-FY-2003"
$ws.Range("A3").Value = "NVDA-FY-2004"
$ws.Range("A4").Value = "NVDA-FY-2005"
$ws.Range("A5").Value = "NVDA-FY-2006"
$ws.Range("A6").Value = "NVDA-FY-2007"
$ws.Range("A7").Value = "NVDA-FY-2008"
$ws.Range("A8").Value = "NVDA-FY-2009"
$ws.Range("A9").Value = "NVDA-FY-2010"
$ws.Range("A10").Value = "NVDA-FY-2011"
$ws.Range("A11").Value = "NVDA-FY-2012"
$ws.Range("A12").Value = "NVDA-FY-2013"
$ws.Range("A13").Value = "NVDA-FY-2014"
$ws.Range("A14").Value = "NVDA-FY-2015"
$ws.Range("A15").Value = "NVDA-FY-2016"
$ws.Range("A16").Value = "NVDA-FY-2017"
$ws.Range("A17").Value = "NVDA-FY-2018"
$ws.Range("A18").Value = "NVDA-FY-2019"
$ws.Range("A19").Value = "NVDA-FY-2020"
$ws.Range("A20").Value = "NVDA-FY-2021"
$ws.Range("A21").Value = "NVDA-FY-2022"

# Numeric data (columns B/C/D) -- FY-2003 row has no data, matching the original pattern for the earliest year
$ws.Range("B3").Value = 2.3125
$ws.Range("C3").Value = 0.7774999737739563
$ws.Range("D3").Value = 1.547145161777735
$ws.Range("B4").Value = 2.27916693687439
$ws.Range("C4").Value = 0.7749999761581421
$ws.Range("D4").Value = 1.60138114111357
$ws.Range("B5").Value = 3.896667003631592
$ws.Range("C5").Value = 1.743332982063293
$ws.Range("D5").Value = 2.505431586527729
$ws.Range("B6").Value = 6.493332862854004
$ws.Range("C6").Value = 2.861666917800903
$ws.Range("D6").Value = 4.666273349761963
$ws.Range("B7").Value = 9.917499542236328
$ws.Range("C7").Value = 4.673333168029785
$ws.Range("D7").Value = 7.001689627414613
$ws.Range("B8").Value = 6.897500038146973
$ws.Range("C8").Value = 1.4375
$ws.Range("D8").Value = 3.616511867922756
$ws.Range("B9").Value = 4.739999771118164
$ws.Range("C9").Value = 1.802500009536743
$ws.Range("D9").Value = 3.143864534290663
$ws.Range("B10").Value = 6.262499809265137
$ws.Range("C10").Value = 2.162499904632568
$ws.Range("D10").Value = 3.469593245831747
$ws.Range("B11").Value = 6.542500019073486
$ws.Range("C11").Value = 2.867500066757202
$ws.Range("D11").Value = 4.087320708658591
$ws.Range("B12").Value = 4.224999904632568
$ws.Range("C12").Value = 2.787499904632568
$ws.Range("D12").Value = 3.341757016967099
$ws.Range("B13").Value = 4.110000133514404
$ws.Range("C13").Value = 3.009999990463257
$ws.Range("D13").Value = 3.605458164595038
$ws.Range("B14").Value = 5.3125
$ws.Range("C14").Value = 3.829999923706055
$ws.Range("D14").Value = 4.699183270275831
$ws.Range("B15").Value = 8.484999656677246
$ws.Range("C15").Value = 4.735000133514404
$ws.Range("D15").Value = 6.087988281622529
$ws.Range("B16").Value = 29.98250007629395
$ws.Range("C16").Value = 6.1875
$ws.Range("D16").Value = 14.82492032374044
$ws.Range("B17").Value = 60.83499908447266
$ws.Range("C17").Value = 23.79249954223633
$ws.Range("D17").Value = 39.58842625560988
$ws.Range("B18").Value = 73.19000244140625
$ws.Range("C18").Value = 31.11499977111816
$ws.Range("D18").Value = 56.77958995819092
$ws.Range("B19").Value = 64.875
$ws.Range("C19").Value = 32.75
$ws.Range("D19").Value = 45.24611545273982
$ws.Range("B20").Value = 147.2675018310547
$ws.Range("C20").Value = 45.16999816894531
$ws.Range("D20").Value = 103.7769730389118
$ws.Range("B21").Value = 346.4700012207031
$ws.Range("C21").Value = 115.6650009155273
$ws.Range("D21").Value = 204.6822222149561

# Column A keeps the same bold/bordered "index" style as the header column (copy from A1)
$ws.Range("A1").Copy()
$ws.Range("A2:A21").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A1").Select()
